$wb = $excel.ActiveWorkbook

# Sheet ALC, row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 19800
$ws.Range("I43").Value = 60000
$ws.Range("K43").Value = 60000
$ws.Range("M43").Value = -59931

# Sheet ALC, row 54
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 15000
$ws.Range("I54").Value = 15000
$ws.Range("K54").Value = 15000
$ws.Range("M54").Value = -14514

# Sheet ALC, row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 90479500
$ws.Range("J62").Value = 4473.75
$ws.Range("L62").Value = 4473.75
$ws.Range("N62").Value = -5721.75

# Sheet ALC, row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 90479500
$ws.Range("J65").Value = 4473.75
$ws.Range("L65").Value = 22368.75
$ws.Range("N65").Value = -28608.75

# Sheet ALC, row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 12505718
$ws.Range("I76").Value = 12505718
$ws.Range("K76").Value = 12505718
$ws.Range("M76").Value = -12505403

# Sheet ALC, row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 12505718
$ws.Range("I79").Value = 12505718
$ws.Range("K79").Value = 12505718
$ws.Range("M79").Value = -12504626

# Sheet ALC, row 92
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 3199
$ws.Range("I92").Value = 3199
$ws.Range("K92").Value = 3199
$ws.Range("M92").Value = -1951

# Sheet ALC, row 100
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 22078.676
$ws.Range("I100").Value = 63709.375
$ws.Range("J100").Value = 9269.23
$ws.Range("K100").Value = 63709.375
$ws.Range("L100").Value = 9269.23
$ws.Range("M100").Value = -63168.375
$ws.Range("N100").Value = -10351.23

# Sheet ALC, row 101
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1537.0769
$ws.Range("J101").Value = 2446.25
$ws.Range("L101").Value = 7338.75
$ws.Range("N101").Value = -10582.75

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3161.1667
$ws.Range("I106").Value = 1630.3636
$ws.Range("K106").Value = 1630.3636
$ws.Range("M106").Value = -999.3635999999999

# Sheet ALC, row 115
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 11111710
$ws.Range("I115").Value = 12346244
$ws.Range("K115").Value = 37038732
$ws.Range("M115").Value = -37037165

# Sheet ALC, row 121
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 1334.6666
$ws.Range("J121").Value = 1334.6666
$ws.Range("L121").Value = 4003.9998
$ws.Range("N121").Value = -7497.9998

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2113
$ws.Range("I141").Value = 1391.25
$ws.Range("K141").Value = 4173.75
$ws.Range("M141").Value = 1006.25

# Sheet ARM, row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46004548
$ws.Range("I2").Value = 46004548
$ws.Range("K2").Value = 46004548
$ws.Range("M2").Value = -46004435

# Sheet ARM, row 21
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 11540
$ws.Range("I21").Value = 475
$ws.Range("J21").Value = 22605
$ws.Range("K21").Value = 475
$ws.Range("L21").Value = 22605
$ws.Range("M21").Value = -101
$ws.Range("N21").Value = -23353

# Sheet ARM, row 106
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 83456.664
$ws.Range("J106").Value = 83456.664
$ws.Range("L106").Value = 83456.664
$ws.Range("N106").Value = -85980.664

# Sheet ARM, row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 46004548
$ws.Range("I116").Value = 46004548
$ws.Range("K116").Value = 46004548
$ws.Range("M116").Value = -46002254

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6605.7144
$ws.Range("I132").Value = 4939.6
$ws.Range("K132").Value = 14818.8
$ws.Range("M132").Value = -12288.8

# Sheet BSM, row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 46004548
$ws.Range("I3").Value = 46004548
$ws.Range("K3").Value = 46004548
$ws.Range("M3").Value = -46004434

# Sheet BSM, row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4034
$ws.Range("I20").Value = 3995.8
$ws.Range("K20").Value = 3995.8
$ws.Range("M20").Value = -3748.8

# Sheet BSM, row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2735.8386
$ws.Range("I134").Value = 1131.92
$ws.Range("K134").Value = 3395.76
$ws.Range("M134").Value = -860.7600000000002

# Sheet CRP, row 28
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 29255.143
$ws.Range("J28").Value = 29255.143
$ws.Range("L28").Value = 29255.143
$ws.Range("N28").Value = -29745.143

# Sheet CRP, row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 98365.234
$ws.Range("I132").Value = 4622
$ws.Range("K132").Value = 13866
$ws.Range("M132").Value = -11336

# Sheet CUL, row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5688409
$ws.Range("I4").Value = 5318510
$ws.Range("J4").Value = 6428206.5
$ws.Range("K4").Value = 15955530
$ws.Range("L4").Value = 19284619.5
$ws.Range("M4").Value = -15955418
$ws.Range("N4").Value = -19284843.5

# Sheet CUL, row 13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 248.2
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 3000
$ws.Range("N13").Value = -3336

# Sheet CUL, row 36
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 551
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

# Sheet GSM, row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1100
$ws.Range("I107").Value = 1100
$ws.Range("K107").Value = 1100
$ws.Range("M107").Value = 820

# Sheet GSM, row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5509.931
$ws.Range("J122").Value = 2993.3333
$ws.Range("L122").Value = 8979.999899999999
$ws.Range("N122").Value = -13879.9999

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4010.7742
$ws.Range("I132").Value = 2439
$ws.Range("K132").Value = 7317
$ws.Range("M132").Value = -4787

# Sheet LTW, row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1843.4166
$ws.Range("I16").Value = 1101.7273
$ws.Range("K16").Value = 1101.7273
$ws.Range("M16").Value = -931.7273

# Sheet LTW, row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 342.7619
$ws.Range("I55").Value = 320
$ws.Range("J55").Value = 379.75
$ws.Range("K55").Value = 320
$ws.Range("L55").Value = 379.75
$ws.Range("M55").Value = -147
$ws.Range("N55").Value = -725.75

# Sheet LTW, row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3716.2307
$ws.Range("J93").Value = 5321.1
$ws.Range("L93").Value = 5321.1
$ws.Range("N93").Value = -7817.1

# Sheet LTW, row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 50005000
$ws.Range("I100").Value = 83336664
$ws.Range("K100").Value = 83336664
$ws.Range("M100").Value = -83336123

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9943.888999999999
$ws.Range("I132").Value = 8249.333000000001
$ws.Range("K132").Value = 24747.999
$ws.Range("M132").Value = -22217.999

# Sheet LTW, row 133
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 49383.375
$ws.Range("J133").Value = 49383.375
$ws.Range("L133").Value = 49383.375
$ws.Range("N133").Value = -54443.375

# Sheet WVR, row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 51876.25
$ws.Range("J96").Value = 51876.25
$ws.Range("L96").Value = 51876.25
$ws.Range("N96").Value = -54622.25

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5917
$ws.Range("I136").Value = 3334
$ws.Range("J136").Value = 8500
$ws.Range("K136").Value = 10002
$ws.Range("L136").Value = 25500
$ws.Range("M136").Value = -7452
